# Update the "想去人数" (interest count) figures in column F across the
# "展览", "演出" and "全部类型" sheets to reflect the newly generated
# gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 835
$ws1.Range("F4").Value  = 318
$ws1.Range("F6").Value  = 1159
$ws1.Range("F12").Value = 55
$ws1.Range("F14").Value = 905
$ws1.Range("F15").Value = 893
$ws1.Range("F17").Value = 77
$ws1.Range("F20").Value = 802
$ws1.Range("F22").Value = 3118
$ws1.Range("F23").Value = 912
$ws1.Range("F25").Value = 2290
$ws1.Range("F27").Value = 6
$ws1.Range("F28").Value = 3147
$ws1.Range("F29").Value = 644
$ws1.Range("F30").Value = 628
$ws1.Range("F32").Value = 93
$ws1.Range("F35").Value = 141
$ws1.Range("F36").Value = 29
$ws1.Range("F38").Value = 1120
$ws1.Range("F39").Value = 1807
$ws1.Range("F40").Value = 410
$ws1.Range("F42").Value = 561
$ws1.Range("F43").Value = 207
$ws1.Range("F46").Value = 53

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 143
$ws2.Range("F12").Value = 92
$ws2.Range("F15").Value = 1
$ws2.Range("F16").Value = 2

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 835
$ws4.Range("F4").Value  = 318
$ws4.Range("F5").Value  = 1159
$ws4.Range("F9").Value  = 55
$ws4.Range("F10").Value = 905
$ws4.Range("F11").Value = 893
$ws4.Range("F12").Value = 143
$ws4.Range("F14").Value = 77
$ws4.Range("F17").Value = 802
$ws4.Range("F19").Value = 3118
$ws4.Range("F20").Value = 912
$ws4.Range("F23").Value = 2290
$ws4.Range("F25").Value = 3147
$ws4.Range("F26").Value = 644
$ws4.Range("F27").Value = 628
$ws4.Range("F33").Value = 93
$ws4.Range("F34").Value = 92
$ws4.Range("F37").Value = 141
$ws4.Range("F41").Value = 1120
$ws4.Range("F42").Value = 1807
$ws4.Range("F43").Value = 2
$ws4.Range("F44").Value = 410
$ws4.Range("F45").Value = 561
$ws4.Range("F46").Value = 207
$ws4.Range("F49").Value = 53
